# Update "想去人数" (want-to-go count) figures in column F on the
# "展览" and "全部类型" worksheets to reflect the newer scrape snapshot.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value  = 13365
$ws1.Range("F6").Value  = 1011
$ws1.Range("F10").Value = 116
$ws1.Range("F11").Value = 68
$ws1.Range("F13").Value = 28
$ws1.Range("F14").Value = 13349
$ws1.Range("F15").Value = 331
$ws1.Range("F16").Value = 585
$ws1.Range("F17").Value = 8895
$ws1.Range("F19").Value = 7968
$ws1.Range("F20").Value = 242
$ws1.Range("F21").Value = 5
$ws1.Range("F26").Value = 18
$ws1.Range("F27").Value = 1014
$ws1.Range("F32").Value = 156

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value  = 13365
$ws4.Range("F7").Value  = 1011
$ws4.Range("F11").Value = 116
$ws4.Range("F12").Value = 68
$ws4.Range("F14").Value = 28
$ws4.Range("F15").Value = 13349
$ws4.Range("F16").Value = 331
$ws4.Range("F17").Value = 585
$ws4.Range("F18").Value = 8895
$ws4.Range("F20").Value = 7968
$ws4.Range("F21").Value = 242
$ws4.Range("F22").Value = 5
$ws4.Range("F27").Value = 18
$ws4.Range("F28").Value = 1014
$ws4.Range("F35").Value = 156
